$wb = $excel.ActiveWorkbook

# --- Rename sheets and drop Sheet3 ---
$wsNPC = $wb.Worksheets.Item("Sheet1")
$wsItemRaw = $wb.Worksheets.Item("Sheet2")
$wb.Worksheets.Item("Sheet3").Delete()

$wsNPC.Name = "NPC"
$wsItemRaw.Name = "Item"

$wsItem = $wb.Worksheets.Item("Item")

# --- Populate the Item sheet, column by column so new shared strings ---
# --- land in the same order as the authored workbook.                ---
$wsItem.Range("A1").Value = "ID"
$wsItem.Range("A2").Value = "string"
$wsItem.Range("A3").Value = "编号"
$wsItem.Range("A4").Value = "WP001"
$wsItem.Range("A5").Value = "WP002"

$wsItem.Range("B1").Value = "Name"
$wsItem.Range("B2").Value = "string"
$wsItem.Range("B3").Value = "名称"
$wsItem.Range("B4").Value = "倚天剑"
$wsItem.Range("B5").Value = "屠龙刀"

$wsItem.Range("C1").Value = "AssetName"
$wsItem.Range("C2").Value = "string"
$wsItem.Range("C3").Value = "资源编号"
$wsItem.Range("C4").Value = "ICON01"
$wsItem.Range("C5").Value = "ICON02"

# Row heights to match the NPC sheet's look
$wsItem.Range("A1:A5").EntireRow.RowHeight = 16.5

# Column C needs to fit the AssetName / ICON0x text
$wsItem.Columns.Item(3).ColumnWidth = 12

# Copy cell formatting from the NPC sheet so the header/data styling matches
$wsNPC.Range("A1:C3").Copy()
$wsItem.Range("A1:C3").PasteSpecial(-4122)

$wsNPC.Range("A4:C5").Copy()
$wsItem.Range("A4:C5").PasteSpecial(-4122)

# --- Selections / active tab ---
# NPC is no longer the active sheet; its old single-cell selection becomes
# a full used-range selection instead.
$wsNPC.Activate()
$wsNPC.Range("A1:XFD6").Select()

# Item becomes the active sheet with C6 selected (just past the data).
$wsItem.Activate()
$wsItem.Range("C6").Select()
